# Applies the "commit before rebase to gh-pages branch" edit:
#   1. Bump the cached "datetimeFigureOut" footer field from 17.10.20 to
#      22.10.20 everywhere it is defined (slide master + every slide layout).
#   2. Nudge four straight-connector shapes on slide 3 down by 22860 EMU
#      (0.025in / 1.8pt) -- only their vertical offset changes.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text: "17.10.20" -> "22.10.20"
# ---------------------------------------------------------------------
function Update-DatePlaceholder {
    param($shapes)

    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "17.10.20") {
                $tr.Text = "22.10.20"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder($master.Shapes)

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder($layouts.Item($i).Shapes)
}

# ---------------------------------------------------------------------
# 2) Move four connectors on slide 3 down by 22860 EMU (keep X/width/height)
# ---------------------------------------------------------------------
# Shape.Top/.Left are expressed in points (1 pt = 12700 EMU) and are stored
# internally as single-precision floats, so a tiny epsilon is added before
# converting so the value rounds back to exactly the target EMU amount.
function Set-ShapeTopEmu {
    param($shape, [double]$targetEmu)

    $pts = $targetEmu / 12700.0
    for ($k = 0; $k -lt 20; $k++) {
        $shape.Top = $pts
        $actualEmu = [math]::Round($shape.Top * 12700.0)
        if ($actualEmu -ge $targetEmu) { break }
        $pts += 0.00001
    }
}

$slide3 = $p.Slides.Item(3)

Set-ShapeTopEmu ($slide3.Shapes.Item("Gerade Verbindung mit Pfeil 14")) 4093830
Set-ShapeTopEmu ($slide3.Shapes.Item("Gerade Verbindung 45"))           4093830
Set-ShapeTopEmu ($slide3.Shapes.Item("Gerade Verbindung mit Pfeil 85")) 4095383
Set-ShapeTopEmu ($slide3.Shapes.Item("Gerade Verbindung 96"))           4095383

Write-Host "Edit complete."
